$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# The underlying data table gained three new data rows:
#   - testModels/refined_db/4/m410/m410.off  (inserted before the
#     existing m412 row, i.e. at worksheet row 179)
#   - testModels/refined_db/4/m411/m411.off  (inserted right after
#     the new m410 row, i.e. at worksheet row 180)
#   - testModels/refined_db/9/m909/m909.off  (inserted before the
#     existing m911 row, which -- after the first shift -- sits at
#     worksheet row 243)
# Every row from the old 179 through 240 shifts down by two rows,
# and the old rows 241-242 shift down by three rows in total.
# We reproduce that by literally inserting blank rows (which moves
# all existing data/formatting down automatically) and then filling
# in the values for the brand-new rows.
# ---------------------------------------------------------------

# Insert two fresh rows at 179 (pushes old row 179 "m412" down to 181)
$ws.Rows.Item(179).Resize(2).Insert()

# Insert one fresh row at 243 (pushes old row 241 "m911", now at 243
# after the previous insert, down to 244)
$ws.Rows.Item(243).Insert()

# Copy the column-A cell formatting onto the freshly inserted, blank
# A-cells so they keep the same style as every other row in the table.
$ws.Range("A178").Copy() | Out-Null
$ws.Range("A179:A180").PasteSpecial(-4122) | Out-Null
$ws.Range("A178").Copy() | Out-Null
$ws.Range("A243").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Fill in the data for the new row 179 -> m410
# ---------------------------------------------------------------
$ws.Range("B179").Value = 4
$ws.Range("C179").Value = 2000
$ws.Range("D179").Value = 1179
$ws.Range("E179").Value = $true
$ws.Range("F179").Value = $false
$ws.Range("G179").Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.44198, 0.94572]))'
$ws.Range("H179").Value = "testModels/refined_db/4/m410/m410.off"
$ws.Range("I179").Value = 0.4426165134704702
$ws.Range("J179").Value = 0.36472700513595
$ws.Range("K179").Value = $false
$ws.Range("L179").Value = $true

# ---------------------------------------------------------------
# Fill in the data for the new row 180 -> m411
# ---------------------------------------------------------------
$ws.Range("B180").Value = 4
$ws.Range("C180").Value = 2000
$ws.Range("D180").Value = 1009
$ws.Range("E180").Value = $true
$ws.Range("F180").Value = $false
$ws.Range("G180").Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.70357, 0.45049]))'
$ws.Range("H180").Value = "testModels/refined_db/4/m411/m411.off"
$ws.Range("I180").Value = 0.356823088341781
$ws.Range("J180").Value = 0.2742882714080499
$ws.Range("K180").Value = $false
$ws.Range("L180").Value = $false

# ---------------------------------------------------------------
# Fill in the data for the new row 243 -> m909
# ---------------------------------------------------------------
$ws.Range("B243").Value = 9
$ws.Range("C243").Value = 2000
$ws.Range("D243").Value = 1010
$ws.Range("E243").Value = $true
$ws.Range("F243").Value = $false
$ws.Range("G243").Value = '(TrackedArray([0.025  , 0.02623, 0.025  ]), TrackedArray([0.70357, 0.44378, 0.975  ]))'
$ws.Range("H243").Value = "testModels/refined_db/9/m909/m909.off"
$ws.Range("I243").Value = 0.1286547542692651
$ws.Range("J243").Value = 0.2691700682120299
$ws.Range("K243").Value = $false
$ws.Range("L243").Value = $false

# ---------------------------------------------------------------
# Column A is a simple running index (row number - 2). Renumber the
# whole column below the header so it stays consistent after the
# inserts above.
# ---------------------------------------------------------------
for ($r = 2; $r -le 245; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------
# Refresh the sheet dimension to match the new extent.
# ---------------------------------------------------------------
$ws.UsedRange | Out-Null
